$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 6739
$ws1.Range("F9").Value = 4577
$ws1.Range("F10").Value = 6794
$ws1.Range("F13").Value = 1395
$ws1.Range("F14").Value = 805
$ws1.Range("F15").Value = 117
$ws1.Range("F20").Value = 130
$ws1.Range("F24").Value = 1060
$ws1.Range("F43").Value = 524
$ws1.Range("F45").Value = 111

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F33").Value = 582

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F9").Value = 1843

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 6739
$ws4.Range("F16").Value = 4577
$ws4.Range("F18").Value = 6794
$ws4.Range("F20").Value = 1395
$ws4.Range("F22").Value = 805
$ws4.Range("F23").Value = 117
$ws4.Range("F27").Value = 130
$ws4.Range("F29").Value = 1060
$ws4.Range("F41").Value = 582
$ws4.Range("F46").Value = 524
$ws4.Range("F49").Value = 111
